$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.514.30'
$ws.Range("E2").Value = '  +5.48%  '
$ws.Range("D3").Value = '1.725.90'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.71'
$ws.Range("E5").Value = '  +3.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5352'
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2668'
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06592'
$ws.Range("E9").Value = '  +4.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.63'
$ws.Range("E10").Value = '  +6.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07699'
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.606'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '1.726.33'
$ws.Range("E13").Value = '  +4.79%  '
$ws.Range("D14").Value = '1.963.15'
$ws.Range("E14").Value = '  +4.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5820'
$ws.Range("E15").Value = '  +4.49%  '
$ws.Range("D16").Value = '0.0₅8272'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.84'
$ws.Range("E17").Value = '  +4.25%  '
$ws.Range("D18").Value = '27.523.88'
$ws.Range("E18").Value = '  +5.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.58'
$ws.Range("E19").Value = '  +13.19%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.726'
$ws.Range("E21").Value = '  +2.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.60'
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.086'
$ws.Range("E23").Value = '  +3.01%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.05'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.757'
$ws.Range("E26").Value = '  +16.25%  '
$ws.Range("E27").Value = '  +4.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.401'
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.54'
$ws.Range("E29").Value = '  +4.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05495'
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.304'
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.563'
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.443'
$ws.Range("E33").Value = '  +3.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.656'
$ws.Range("E34").Value = '  +6.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.859'
$ws.Range("E35").Value = '  +2.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9651'
$ws.Range("E36").Value = '  +2.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.425'
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5980'
$ws.Range("E38").Value = '  +6.98%  '
$ws.Range("E39").Value = '  +4.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.896'
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("D41").Value = '1.055.68'
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8527'
$ws.Range("E42").Value = '  +3.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.004'
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.32'
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").Value = '1.870.00'
$ws.Range("E45").Value = '  +4.72%  '
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '58.93'
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4472'
$ws.Range("E48").Value = '  +3.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.185'
$ws.Range("E49").Value = '  +3.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.003'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05243'
$ws.Range("E51").Value = '  +2.61%  '
